$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 2).Value = 0.5867794862083637
$ws.Cells.Item(4, 3).Value = 0.588
$ws.Cells.Item(4, 4).Value = 0.5911496908674081
$ws.Cells.Item(4, 5).Value = 0.59
$ws.Cells.Item(4, 6).Value = 0.6755282899659024
$ws.Cells.Item(4, 7).Value = 0.6990000000000001
$ws.Cells.Item(4, 8).Value = 0.6588925889241295
$ws.Cells.Item(4, 9).Value = 0.664
$ws.Cells.Item(4, 10).Value = 0.4888536848760653
$ws.Cells.Item(4, 11).Value = 0.4970000000000001
$ws.Cells.Item(4, 12).Value = 0.4853343332709336
$ws.Cells.Item(4, 13).Value = 0.483
$ws.Cells.Item(4, 14).Value = 0.6580453174800088
$ws.Cells.Item(4, 15).Value = 0.6859999999999999
$ws.Cells.Item(4, 16).Value = 0.6385261638611831
$ws.Cells.Item(4, 17).Value = 0.647
$ws.Cells.Item(4, 18).Value = 0.5910910638846751
$ws.Cells.Item(4, 19).Value = 0.5940000000000001
$ws.Cells.Item(4, 20).Value = 0.5930100494528215
$ws.Cells.Item(4, 21).Value = 0.5924999999999999
$ws.Cells.Item(4, 22).Value = 0.6714857468567652
$ws.Cells.Item(4, 23).Value = 0.696
$ws.Cells.Item(4, 24).Value = 0.654044396947697
$ws.Cells.Item(4, 25).Value = 0.6599999999999999
$ws.Cells.Item(4, 26).Value = 0.6623944431189186
$ws.Cells.Item(4, 27).Value = 0.6910000000000001
$ws.Cells.Item(4, 28).Value = 0.6416551457567504
$ws.Cells.Item(4, 29).Value = 0.6505

$ws.Cells.Item(5, 2).Value = 0.6524668262824616
$ws.Cells.Item(5, 3).Value = 0.8470000000000001
$ws.Cells.Item(5, 4).Value = 0.5342609065801748
$ws.Cells.Item(5, 5).Value = 0.5525
$ws.Cells.Item(5, 6).Value = 0.6623177022921712
$ws.Cells.Item(5, 7).Value = 0.849
$ws.Cells.Item(5, 8).Value = 0.5511598390280821
$ws.Cells.Item(5, 9).Value = 0.5740000000000001
$ws.Cells.Item(5, 10).Value = 0.6687938021972293
$ws.Cells.Item(5, 11).Value = 0.946
$ws.Cells.Item(5, 12).Value = 0.5182224495543275
$ws.Cells.Item(5, 13).Value = 0.532
$ws.Cells.Item(5, 14).Value = 0.6540053142392115
$ws.Cells.Item(5, 15).Value = 0.852
$ws.Cells.Item(5, 16).Value = 0.5355237034507258
$ws.Cells.Item(5, 17).Value = 0.5534999999999999
$ws.Cells.Item(5, 18).Value = 0.6561070483877145
$ws.Cells.Item(5, 19).Value = 0.86
$ws.Cells.Item(5, 20).Value = 0.5339164801115952
$ws.Cells.Item(5, 21).Value = 0.553
$ws.Cells.Item(5, 22).Value = 0.6702176239762629
$ws.Cells.Item(5, 23).Value = 0.885
$ws.Cells.Item(5, 24).Value = 0.5450095930567206
$ws.Cells.Item(5, 25).Value = 0.57
$ws.Cells.Item(5, 26).Value = 0.6553253939830184
$ws.Cells.Item(5, 27).Value = 0.861
$ws.Cells.Item(5, 28).Value = 0.5334430843523295
$ws.Cells.Item(5, 29).Value = 0.5515000000000001

$ws.Cells.Item(6, 2).Value = 0.5743952573160535
$ws.Cells.Item(6, 3).Value = 0.5620000000000001
$ws.Cells.Item(6, 4).Value = 0.5997305832069287
$ws.Cells.Item(6, 5).Value = 0.6010000000000001
$ws.Cells.Item(6, 6).Value = 0.6935710155645171
$ws.Cells.Item(6, 7).Value = 0.6960000000000001
$ws.Cells.Item(6, 8).Value = 0.7070314191332031
$ws.Cells.Item(6, 9).Value = 0.6945
$ws.Cells.Item(6, 10).Value = 0.5193020752264974
$ws.Cells.Item(6, 11).Value = 0.5240000000000001
$ws.Cells.Item(6, 12).Value = 0.5244858522170122
$ws.Cells.Item(6, 13).Value = 0.5225000000000001
$ws.Cells.Item(6, 14).Value = 0.6676778392119738
$ws.Cells.Item(6, 15).Value = 0.6499999999999999
$ws.Cells.Item(6, 16).Value = 0.710395504594944
$ws.Cells.Item(6, 17).Value = 0.6899999999999999
$ws.Cells.Item(6, 18).Value = 0.5829022968259514
$ws.Cells.Item(6, 19).Value = 0.5700000000000001
$ws.Cells.Item(6, 20).Value = 0.6113618683618081
$ws.Cells.Item(6, 21).Value = 0.6105
$ws.Cells.Item(6, 22).Value = 0.6853297188491376
$ws.Cells.Item(6, 23).Value = 0.6969999999999998
$ws.Cells.Item(6, 24).Value = 0.6863753167185024
$ws.Cells.Item(6, 25).Value = 0.6830000000000001
$ws.Cells.Item(6, 26).Value = 0.6795534778210511
$ws.Cells.Item(6, 27).Value = 0.6599999999999999
$ws.Cells.Item(6, 28).Value = 0.7226943761333404
$ws.Cells.Item(6, 29).Value = 0.7

Write-Host "Updated classifier results"